$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 gets rebuilt: two new cards (blog #139 + the Product Hunt
# launch "signin" card) are inserted at the front, pushing every existing
# card two columns to the right, with J/K (ser:138 vs the sign-up card)
# swapping back to their original relative order, and the old C9 value
# ("ser: 137") sliding down to the new L9 / old L9 (footer) sliding to M9.
#
# Read every old value into variables FIRST (before any writes) so the
# shared-string table never transiently loses the only reference to a
# value we still need to relocate.
$oldB = $ws.Range("B9").Value()
$oldC = $ws.Range("C9").Value()
$oldD = $ws.Range("D9").Value()
$oldE = $ws.Range("E9").Value()
$oldF = $ws.Range("F9").Value()
$oldG = $ws.Range("G9").Value()
$oldH = $ws.Range("H9").Value()
$oldJ = $ws.Range("J9").Value()
$oldK = $ws.Range("K9").Value()
$oldL = $ws.Range("L9").Value()

# Write back-to-front so a cell is never overwritten before its old value
# has already been copied to its new home.
$ws.Range("M9").Value = $oldL
$ws.Range("L9").Value = $oldC
$ws.Range("K9").Value = $oldJ
$ws.Range("J9").Value = $oldK
$ws.Range("I9").Value = $oldH
$ws.Range("H9").Value = $oldG
$ws.Range("G9").Value = $oldF
$ws.Range("F9").Value = $oldE
$ws.Range("E9").Value = $oldD
$ws.Range("D9").Value = $oldB

# New content for B9 / C9 (brand-new shared strings). The "ser: 139" blog
# card (shared-string index 33) must land in the table before the Product
# Hunt signin card (index 34), so C9 — which ends up holding the blog card
# — is written before B9, matching the author's shared-string order.
$ws.Range("C9").Value = "type: blog`nwidth: 2`nheight: 1`nser: 139"
$ws.Range("B9").Value = "type: signin`nwidth: 4`nheight: 1`nh3: Zakatlists is Launching on Product Hunt - 27 May 2020`np.w-m-50: I have been writing blogs for nearly 140 days this year. If these blogs or our techshek conferences have impacted your life in some way, please support me on the launch day. I will remind you on phone call or email you if you feel you will forget. `nbutton.default: Set a reminder*goto(`"https://docs.google.com/forms/d/e/1FAIpQLScRWGicOlVW-RpquUYvHD769v45XACKZydnpTJcJVWxzHW0jg/viewform?usp=sf_link`")`nbutton.default: What is Product Hunt*goto(`"https://www.producthunt.com/about`")`nsvg: /icons/producthunt.svg"

# Make sure the wrap-text style used by the rest of the data rows carries
# onto the new M column cell too.
$ws.Range("M9").WrapText = $true

# Row grows to the (already-present-elsewhere) max auto height.
$ws.Range("A9:M9").EntireRow.RowHeight = 409.6

# New column M needs the same kind of width definition as K/L.
$ws.Columns.Item(13).ColumnWidth = 35.6640625

# Selection / view bookkeeping to mirror the author's final state.
$ws.Range("D9").Select()
$excel.ActiveWindow.ScrollRow = 1
